$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Target cluster (column D) text values
$ws.Range("D2").Value = "FAPs"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("D5").Value = "FAPs"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("D8").Value = "FAPs"
$ws.Range("D9").Value = "MuSCs"
$ws.Range("D10").Value = "Resolving-Mac"
$ws.Range("D11").Value = "FAPs"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("D13").Value = "Resolving-Mac"

# Update numeric columns (G-T) with recomputed TPM-based values
# Row 2
$ws.Range("G2").Value = 54.53585066666667
$ws.Range("H2").Value = 163.607552
$ws.Range("I2").Value = 0.3031388658437607
$ws.Range("J2").Value = 0.3031388658437607
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.987076
$ws.Range("N2").Value = 11.961228
$ws.Range("O2").Value = 0.2813308272685638
$ws.Range("P2").Value = 0.2813308272685638
$ws.Range("Q2").Value = 217.4385813326507
$ws.Range("R2").Value = 1956.947231993856
$ws.Range("S2").Value = 0.08528230790507936
$ws.Range("T2").Value = 0.08528230790507937
# Row 3
$ws.Range("G3").Value = 54.53585066666667
$ws.Range("H3").Value = 163.607552
$ws.Range("I3").Value = 0.3031388658437607
$ws.Range("J3").Value = 0.3031388658437607
$ws.Range("M3").Value = 10.131229
$ws.Range("N3").Value = 30.393687
$ws.Range("O3").Value = 0.7148664925918803
$ws.Range("P3").Value = 0.7148664925918804
$ws.Range("Q3").Value = 552.5151918138026
$ws.Range("R3").Value = 4972.636726324224
$ws.Range("S3").Value = 0.2167038177940097
$ws.Range("T3").Value = 0.2167038177940098
# Row 4
$ws.Range("G4").Value = 54.53585066666667
$ws.Range("H4").Value = 163.607552
$ws.Range("I4").Value = 0.3031388658437607
$ws.Range("J4").Value = 0.3031388658437607
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05389233333333333
$ws.Range("N4").Value = 0.161677
$ws.Range("O4").Value = 0.00380268013955587
$ws.Range("P4").Value = 0.00380268013955587
$ws.Range("Q4").Value = 2.939064242744889
$ws.Range("R4").Value = 26.451578184704
$ws.Range("S4").Value = 0.00115274014467156
$ws.Range("T4").Value = 0.00115274014467156
# Row 5
$ws.Range("I5").Value = 0.1026363515063155
$ws.Range("J5").Value = 0.1026363515063155
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.987076
$ws.Range("N5").Value = 11.961228
$ws.Range("O5").Value = 0.2813308272685638
$ws.Range("P5").Value = 0.2813308272685638
$ws.Range("Q5").Value = 73.62006386932534
$ws.Range("R5").Value = 662.5805748239279
$ws.Range("S5").Value = 0.02887476967709884
$ws.Range("T5").Value = 0.02887476967709884
# Row 6
$ws.Range("I6").Value = 0.1026363515063155
$ws.Range("J6").Value = 0.1026363515063155
$ws.Range("M6").Value = 10.131229
$ws.Range("N6").Value = 30.393687
$ws.Range("O6").Value = 0.7148664925918803
$ws.Range("P6").Value = 0.7148664925918804
$ws.Range("Q6").Value = 187.0698542126513
$ws.Range("R6").Value = 1683.628687913862
$ws.Range("S6").Value = 0.0733712886137471
$ws.Range("T6").Value = 0.07337128861374713
# Row 7
$ws.Range("I7").Value = 0.1026363515063155
$ws.Range("J7").Value = 0.1026363515063155
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.05389233333333333
$ws.Range("N7").Value = 0.161677
$ws.Range("O7").Value = 0.00380268013955587
$ws.Range("P7").Value = 0.00380268013955587
$ws.Range("Q7").Value = 0.9951044379557776
$ws.Range("R7").Value = 8.955939941601999
$ws.Range("S7").Value = 0.0003902932154695411
$ws.Range("T7").Value = 0.0003902932154695412
# Row 8
$ws.Range("G8").Value = 12.55635966666667
$ws.Range("H8").Value = 37.669079
$ws.Range("I8").Value = 0.06979483370938171
$ws.Range("J8").Value = 0.06979483370938172
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.987076
$ws.Range("N8").Value = 11.961228
$ws.Range("O8").Value = 0.2813308272685638
$ws.Range("P8").Value = 0.2813308272685638
$ws.Range("Q8").Value = 50.06316027433466
$ws.Range("R8").Value = 450.568442469012
$ws.Range("S8").Value = 0.0196354383065322
$ws.Range("T8").Value = 0.0196354383065322
# Row 9
$ws.Range("G9").Value = 12.55635966666667
$ws.Range("H9").Value = 37.669079
$ws.Range("I9").Value = 0.06979483370938171
$ws.Range("J9").Value = 0.06979483370938172
$ws.Range("M9").Value = 10.131229
$ws.Range("N9").Value = 30.393687
$ws.Range("O9").Value = 0.7148664925918803
$ws.Range("P9").Value = 0.7148664925918804
$ws.Range("Q9").Value = 127.2113551893636
$ws.Range("R9").Value = 1144.902196704273
$ws.Range("S9").Value = 0.04989398797485924
$ws.Range("T9").Value = 0.04989398797485925
# Row 10
$ws.Range("G10").Value = 12.55635966666667
$ws.Range("H10").Value = 37.669079
$ws.Range("I10").Value = 0.06979483370938171
$ws.Range("J10").Value = 0.06979483370938172
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.05389233333333333
$ws.Range("N10").Value = 0.161677
$ws.Range("O10").Value = 0.00380268013955587
$ws.Range("P10").Value = 0.00380268013955587
$ws.Range("Q10").Value = 0.6766915206092221
$ws.Range("R10").Value = 6.090223685482999
$ws.Range("S10").Value = 0.0002654074279902704
$ws.Range("T10").Value = 0.0002654074279902705
# Row 11
$ws.Range("G11").Value = 94.34696966666667
$ws.Range("H11").Value = 283.040909
$ws.Range("I11").Value = 0.524429948940542
$ws.Range("J11").Value = 0.5244299489405421
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.987076
$ws.Range("N11").Value = 11.961228
$ws.Range("O11").Value = 0.2813308272685638
$ws.Range("P11").Value = 0.2813308272685638
$ws.Range("Q11").Value = 376.1685384306947
$ws.Range("R11").Value = 3385.516845876252
$ws.Range("S11").Value = 0.1475383113798533
$ws.Range("T11").Value = 0.1475383113798534
# Row 12
$ws.Range("G12").Value = 94.34696966666667
$ws.Range("H12").Value = 283.040909
$ws.Range("I12").Value = 0.524429948940542
$ws.Range("J12").Value = 0.5244299489405421
$ws.Range("M12").Value = 10.131229
$ws.Range("N12").Value = 30.393687
$ws.Range("O12").Value = 0.7148664925918803
$ws.Range("P12").Value = 0.7148664925918804
$ws.Range("Q12").Value = 955.8507551490536
$ws.Range("R12").Value = 8602.656796341484
$ws.Range("S12").Value = 0.3748973982092641
$ws.Range("T12").Value = 0.3748973982092642
# Row 13
$ws.Range("G13").Value = 94.34696966666667
$ws.Range("H13").Value = 283.040909
$ws.Range("I13").Value = 0.524429948940542
$ws.Range("J13").Value = 0.5244299489405421
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.05389233333333333
$ws.Range("N13").Value = 0.161677
$ws.Range("O13").Value = 0.00380268013955587
$ws.Range("P13").Value = 0.00380268013955587
$ws.Range("Q13").Value = 5.084578338265889
$ws.Range("R13").Value = 45.76120504439299
$ws.Range("S13").Value = 0.001994239351424498
$ws.Range("T13").Value = 0.001994239351424499
